$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all.inputs")

$ws.Cells.Item(1, 1).Value = '$all.params'
$ws.Cells.Item(2, 1).Value = '       r0.initial latent.period illness.length.given.nonhosp prop.hospitalized'
$ws.Cells.Item(3, 1).Value = '    1:   2.574271             4                            3       0.043527035'
$ws.Cells.Item(4, 1).Value = '    2:   1.793005             0                            4       0.084637599'
$ws.Cells.Item(5, 1).Value = '    3:   3.996326             2                            4       0.037620719'
$ws.Cells.Item(6, 1).Value = '    4:   2.934332             1                            5       0.016961280'
$ws.Cells.Item(7, 1).Value = '    5:   6.120791             2                            8       0.034717478'
$ws.Cells.Item(8, 1).Value = '   ---                                                                        '
$ws.Cells.Item(9, 1).Value = ' 9996:   3.424709             1                            6       0.025086411'
$ws.Cells.Item(10, 1).Value = ' 9997:   2.641475             5                            6       0.052582408'
$ws.Cells.Item(11, 1).Value = ' 9998:   4.629994             5                            4       0.052749284'
$ws.Cells.Item(12, 1).Value = ' 9999:   4.184294             6                            4       0.057032881'
$ws.Cells.Item(13, 1).Value = '10000:   3.014694             3                            3       0.002760913'
$ws.Cells.Item(14, 1).Value = '       hosp.length.of.stay  prop.icu prop.vent intervention1.date'
$ws.Cells.Item(15, 1).Value = '    1:                  18 0.4522037 0.8146784         2020-03-07'
$ws.Cells.Item(16, 1).Value = '    2:                  23 0.6040443 0.9049126         2020-03-07'
$ws.Cells.Item(17, 1).Value = '    3:                  17 0.4785787 0.7982800         2020-03-07'
$ws.Cells.Item(18, 1).Value = '    4:                   9 0.4221481 0.8458649         2020-03-07'
$ws.Cells.Item(19, 1).Value = '    5:                  16 0.4894000 0.8549066         2020-03-07'
$ws.Cells.Item(20, 1).Value = '   ---                                                           '
$ws.Cells.Item(21, 1).Value = ' 9996:                  17 0.4381462 0.8991000         2020-03-07'
$ws.Cells.Item(22, 1).Value = ' 9997:                   8 0.5170023 0.8220893         2020-03-07'
$ws.Cells.Item(23, 1).Value = ' 9998:                  19 0.4848476 0.8123715         2020-03-07'
$ws.Cells.Item(24, 1).Value = ' 9999:                   4 0.4802124 0.7861673         2020-03-07'
$ws.Cells.Item(25, 1).Value = '10000:                  18 0.4557356 0.9204351         2020-03-07'
$ws.Cells.Item(26, 1).Value = '       intervention1.multiplier intervention1.smooth.days intervention2.date'
$ws.Cells.Item(27, 1).Value = '    1:                0.3933902                         3         2020-03-17'
$ws.Cells.Item(28, 1).Value = '    2:                0.5162912                         5         2020-03-17'
$ws.Cells.Item(29, 1).Value = '    3:                0.6272024                         4         2020-03-17'
$ws.Cells.Item(30, 1).Value = '    4:                0.6729379                         7         2020-03-17'
$ws.Cells.Item(31, 1).Value = '    5:                0.3846296                        10         2020-03-17'
$ws.Cells.Item(32, 1).Value = '   ---                                                                      '
$ws.Cells.Item(33, 1).Value = ' 9996:                0.6843684                         4         2020-03-17'
$ws.Cells.Item(34, 1).Value = ' 9997:                0.3081265                        11         2020-03-17'
$ws.Cells.Item(35, 1).Value = ' 9998:                0.5053716                        11         2020-03-17'
$ws.Cells.Item(36, 1).Value = ' 9999:                0.4312530                         2         2020-03-17'
$ws.Cells.Item(37, 1).Value = '10000:                0.2939671                         8         2020-03-17'
$ws.Cells.Item(38, 1).Value = '       intervention2.multiplier intervention2.smooth.days intervention3.date'
$ws.Cells.Item(39, 1).Value = '    1:                0.4055601                         9         2020-06-01'
$ws.Cells.Item(40, 1).Value = '    2:                0.5106207                         7         2020-06-01'
$ws.Cells.Item(41, 1).Value = '    3:                0.2840808                         3         2020-06-01'
$ws.Cells.Item(42, 1).Value = '    4:                0.6510927                         4         2020-06-01'
$ws.Cells.Item(43, 1).Value = '    5:                0.7780911                         6         2020-06-01'
$ws.Cells.Item(44, 1).Value = '   ---                                                                      '
$ws.Cells.Item(45, 1).Value = ' 9996:                0.2580803                         6         2020-06-01'
$ws.Cells.Item(46, 1).Value = ' 9997:                0.4810451                         9         2020-06-01'
$ws.Cells.Item(47, 1).Value = ' 9998:                0.6432268                         4         2020-06-01'
$ws.Cells.Item(48, 1).Value = ' 9999:                0.3043401                         6         2020-06-01'
$ws.Cells.Item(49, 1).Value = '10000:                1.2236223                         6         2020-06-01'
$ws.Cells.Item(50, 1).Value = '       intervention3.multiplier intervention3.smooth.days'
$ws.Cells.Item(51, 1).Value = '    1:                 1.788350                         7'
$ws.Cells.Item(52, 1).Value = '    2:                 1.284699                         8'
$ws.Cells.Item(53, 1).Value = '    3:                 1.454884                         4'
$ws.Cells.Item(54, 1).Value = '    4:                 2.368237                        12'
$ws.Cells.Item(55, 1).Value = '    5:                 1.419332                         6'
$ws.Cells.Item(56, 1).Value = '   ---                                                   '
$ws.Cells.Item(57, 1).Value = ' 9996:                 1.567420                         3'
$ws.Cells.Item(58, 1).Value = ' 9997:                 1.440978                         7'
$ws.Cells.Item(59, 1).Value = ' 9998:                 1.362914                         3'
$ws.Cells.Item(60, 1).Value = ' 9999:                 2.325578                         7'
$ws.Cells.Item(61, 1).Value = '10000:                 1.670447                         3'
$ws.Cells.Item(62, 1).Value = '       patients.in.hosp.are.infectious use.hosp.rate exposed.to.hospital'
$ws.Cells.Item(63, 1).Value = '    1:                           FALSE         FALSE                  10'
$ws.Cells.Item(64, 1).Value = '    2:                           FALSE         FALSE                   8'
$ws.Cells.Item(65, 1).Value = '    3:                            TRUE         FALSE                   4'
$ws.Cells.Item(66, 1).Value = '    4:                            TRUE          TRUE                   8'
$ws.Cells.Item(67, 1).Value = '    5:                           FALSE         FALSE                   6'
$ws.Cells.Item(68, 1).Value = '   ---                                                                  '
$ws.Cells.Item(69, 1).Value = ' 9996:                           FALSE         FALSE                   9'
$ws.Cells.Item(70, 1).Value = ' 9997:                            TRUE         FALSE                  11'
$ws.Cells.Item(71, 1).Value = ' 9998:                           FALSE         FALSE                  12'
$ws.Cells.Item(72, 1).Value = ' 9999:                           FALSE         FALSE                  12'
$ws.Cells.Item(73, 1).Value = '10000:                           FALSE          TRUE                   8'
$ws.Cells.Item(74, 1).Value = ''
$ws.Cells.Item(75, 1).Value = '$devlist'
$ws.Cells.Item(76, 1).Value = 'quartz_off_screen '
$ws.Cells.Item(77, 1).Value = '                2 '
$ws.Cells.Item(78, 1).Value = ''
$ws.Cells.Item(79, 1).Value = '$extras'
$ws.Cells.Item(80, 1).Value = '$extras$`Parameters with Distributions`'
$ws.Cells.Item(81, 1).Value = '                      internal.name'
$ws.Cells.Item(82, 1).Value = ' 1:                   weight.labels'
$ws.Cells.Item(83, 1).Value = ' 2:               parameter.weights'
$ws.Cells.Item(84, 1).Value = ' 3:                      r0.initial'
$ws.Cells.Item(85, 1).Value = ' 4:                   latent.period'
$ws.Cells.Item(86, 1).Value = ' 5:    illness.length.given.nonhosp'
$ws.Cells.Item(87, 1).Value = ' 6:          infectious.to.hospital'
$ws.Cells.Item(88, 1).Value = ' 7:               prop.hospitalized'
$ws.Cells.Item(89, 1).Value = ' 8:             hosp.length.of.stay'
$ws.Cells.Item(90, 1).Value = ' 9:                        prop.icu'
$ws.Cells.Item(91, 1).Value = '10:                       prop.vent'
$ws.Cells.Item(92, 1).Value = '11:              intervention1.date'
$ws.Cells.Item(93, 1).Value = '12:        intervention1.multiplier'
$ws.Cells.Item(94, 1).Value = '13:       intervention1.smooth.days'
$ws.Cells.Item(95, 1).Value = '14:              intervention2.date'
$ws.Cells.Item(96, 1).Value = '15:        intervention2.multiplier'
$ws.Cells.Item(97, 1).Value = '16:       intervention2.smooth.days'
$ws.Cells.Item(98, 1).Value = '17:              intervention3.date'
$ws.Cells.Item(99, 1).Value = '18:        intervention3.multiplier'
$ws.Cells.Item(100, 1).Value = '19:       intervention3.smooth.days'
$ws.Cells.Item(101, 1).Value = '20: patients.in.hosp.are.infectious'
$ws.Cells.Item(102, 1).Value = '21:                   use.hosp.rate'
$ws.Cells.Item(103, 1).Value = '                      internal.name'
$ws.Cells.Item(104, 1).Value = '                                                              external.name'
$ws.Cells.Item(105, 1).Value = ' 1:                                                                    <NA>'
$ws.Cells.Item(106, 1).Value = ' 2:                                                                  Priors'
$ws.Cells.Item(107, 1).Value = ' 3:                       Basic reproductive number R0 before Intervention1'
$ws.Cells.Item(108, 1).Value = ' 4:    Number of Days from Infection to Becoming Infectious (Latent Period)'
$ws.Cells.Item(109, 1).Value = ' 5:                                       Duration of infectiousness (days)'
$ws.Cells.Item(110, 1).Value = ' 6:             Time from onset of infectiousness to hospitalization (days)'
$ws.Cells.Item(111, 1).Value = ' 7:                               Percent of Infected that are Hospitalized'
$ws.Cells.Item(112, 1).Value = ' 8:                                  Average Hospital Length of Stay (Days)'
$ws.Cells.Item(113, 1).Value = ' 9: Percent of Hospitalized COVID-19 Patients That are Currently in the ICU'
$ws.Cells.Item(114, 1).Value = '10:    Percent of COVID-19 Patients in the ICU who are Currently Ventilated'
$ws.Cells.Item(115, 1).Value = '11:                                              Date of first intervention'
$ws.Cells.Item(116, 1).Value = '12:                                                           Re multiplier'
$ws.Cells.Item(117, 1).Value = '13:                                                    Days to reach new Re'
$ws.Cells.Item(118, 1).Value = '14:                                             Date of second intervention'
$ws.Cells.Item(119, 1).Value = '15:                                                           Re multiplier'
$ws.Cells.Item(120, 1).Value = '16:                                                    Days to reach new Re'
$ws.Cells.Item(121, 1).Value = '17:                                              Date of third intervention'
$ws.Cells.Item(122, 1).Value = '18:                                                           Re multiplier'
$ws.Cells.Item(123, 1).Value = '19:                                                    Days to reach new Re'
$ws.Cells.Item(124, 1).Value = '20:                                     Patients in hospital are infectious'
$ws.Cells.Item(125, 1).Value = '21:   Contant rate to hospital (if FALSE, fixed number of days to hospital)'
$ws.Cells.Item(126, 1).Value = '                                                              external.name'
$ws.Cells.Item(127, 1).Value = '           low     midlow                 mid    midhigh       high'
$ws.Cells.Item(128, 1).Value = ' 1:         NA         NA User''s "Best Guess"         NA         NA'
$ws.Cells.Item(129, 1).Value = ' 2:        0.2        0.2                 0.2        0.2        0.2'
$ws.Cells.Item(130, 1).Value = ' 3:        2.5          3                   4          4        4.5'
$ws.Cells.Item(131, 1).Value = ' 4:          0          2                   3          4          5'
$ws.Cells.Item(132, 1).Value = ' 5:          3          4                   5          6          7'
$ws.Cells.Item(133, 1).Value = ' 6:          4          5                   6          7          8'
$ws.Cells.Item(134, 1).Value = ' 7:       0.01       0.02                0.04       0.05       0.06'
$ws.Cells.Item(135, 1).Value = ' 8:          6         10                  14         18         22'
$ws.Cells.Item(136, 1).Value = ' 9:        0.4       0.42                0.45       0.52       0.55'
$ws.Cells.Item(137, 1).Value = '10:        0.8       0.82                0.85       0.87        0.9'
$ws.Cells.Item(138, 1).Value = '11: 2020-03-07 2020-03-07          2020-03-07 2020-03-07 2020-03-07'
$ws.Cells.Item(139, 1).Value = '12:       0.35        0.4                0.45       0.55        0.7'
$ws.Cells.Item(140, 1).Value = '13:          3          5                   7          9         11'
$ws.Cells.Item(141, 1).Value = '14: 2020-03-17 2020-03-17          2020-03-17 2020-03-17 2020-03-17'
$ws.Cells.Item(142, 1).Value = '15:        0.3        0.4                0.45        0.8          1'
$ws.Cells.Item(143, 1).Value = '16:          3          5                   7          9         11'
$ws.Cells.Item(144, 1).Value = '17: 2020-06-01 2020-06-01          2020-06-01 2020-06-01 2020-06-01'
$ws.Cells.Item(145, 1).Value = '18:        1.1        1.2                 1.5        1.7          2'
$ws.Cells.Item(146, 1).Value = '19:          3          5                   7          9         11'
$ws.Cells.Item(147, 1).Value = '20:      FALSE      FALSE                TRUE      FALSE      FALSE'
$ws.Cells.Item(148, 1).Value = '21:      FALSE      FALSE                TRUE      FALSE      FALSE'
$ws.Cells.Item(149, 1).Value = '           low     midlow                 mid    midhigh       high'
$ws.Cells.Item(150, 1).Value = ''
$ws.Cells.Item(151, 1).Value = '$extras$`Model Inputs`'
$ws.Cells.Item(152, 1).Value = '        internal.name            external.name      value'
$ws.Cells.Item(153, 1).Value = '1:   total.population Number of People in Area    1671000'
$ws.Cells.Item(154, 1).Value = '2: start.display.date    Projection Start Date 2020-03-15'
$ws.Cells.Item(155, 1).Value = '3:           end.date      Projection End Date 2020-07-01'
$ws.Cells.Item(156, 1).Value = ''
$ws.Cells.Item(157, 1).Value = '$extras$`Hospitilization Data`'
$ws.Cells.Item(158, 1).Value = '          Date LowerBound UpperBound'
$ws.Cells.Item(159, 1).Value = ' 1: 2020-04-01   51.88321   59.66569'
$ws.Cells.Item(160, 1).Value = ' 2: 2020-04-02   56.86563   65.39548'
$ws.Cells.Item(161, 1).Value = ' 3: 2020-04-03   61.68362   70.93617'
$ws.Cells.Item(162, 1).Value = ' 4: 2020-04-04   66.26620   76.20613'
$ws.Cells.Item(163, 1).Value = ' 5: 2020-04-05   70.60557   81.19641'
$ws.Cells.Item(164, 1).Value = ' 6: 2020-04-06   74.69862   85.90341'
$ws.Cells.Item(165, 1).Value = ' 7: 2020-04-07   78.82360   90.64714'
$ws.Cells.Item(166, 1).Value = ' 8: 2020-04-08   82.50392   94.87951'
$ws.Cells.Item(167, 1).Value = ' 9: 2020-04-09   85.20577   97.98664'
$ws.Cells.Item(168, 1).Value = '10: 2020-04-10   87.19809  100.27781'
$ws.Cells.Item(169, 1).Value = '11: 2020-04-11   88.66911  101.96948'
$ws.Cells.Item(170, 1).Value = '12: 2020-04-12   89.34864  102.75094'
$ws.Cells.Item(171, 1).Value = '13: 2020-04-13   88.25265  101.49055'
$ws.Cells.Item(172, 1).Value = '14: 2020-04-14   86.63990   99.63588'
$ws.Cells.Item(173, 1).Value = '15: 2020-04-15   85.23594   98.02133'
$ws.Cells.Item(174, 1).Value = '16: 2020-04-17   84.00665   96.60765'
$ws.Cells.Item(175, 1).Value = '17: 2020-04-18   83.83743   96.41305'
$ws.Cells.Item(176, 1).Value = '18: 2020-04-19   83.86506   96.44481'
$ws.Cells.Item(177, 1).Value = '19: 2020-04-20   83.89666   96.48116'
$ws.Cells.Item(178, 1).Value = '20: 2020-04-21   83.85270   96.43061'
$ws.Cells.Item(179, 1).Value = '21: 2020-04-22   83.14715   95.61923'
$ws.Cells.Item(180, 1).Value = '22: 2020-04-23   82.40143   94.76164'
$ws.Cells.Item(181, 1).Value = '23: 2020-04-24   82.17643   94.50289'
$ws.Cells.Item(182, 1).Value = '24: 2020-04-25   81.80076   94.07088'
$ws.Cells.Item(183, 1).Value = '25: 2020-04-26   81.04030   93.19635'
$ws.Cells.Item(184, 1).Value = '26: 2020-04-27   80.07035   92.08090'
$ws.Cells.Item(185, 1).Value = '27: 2020-04-28   78.38229   90.13963'
$ws.Cells.Item(186, 1).Value = '28: 2020-04-29   76.64307   88.13953'
$ws.Cells.Item(187, 1).Value = '29: 2020-04-30   75.31391   86.61100'
$ws.Cells.Item(188, 1).Value = '30: 2020-05-01   74.15291   85.27584'
$ws.Cells.Item(189, 1).Value = '31: 2020-05-02   73.32236   84.32072'
$ws.Cells.Item(190, 1).Value = '32: 2020-05-03   72.73400   83.64410'
$ws.Cells.Item(191, 1).Value = '          Date LowerBound UpperBound'
$ws.Cells.Item(192, 1).Value = ''
$ws.Cells.Item(193, 1).Value = '$extras$Internal'
$ws.Cells.Item(194, 1).Value = '                    internal.name             value'
$ws.Cells.Item(195, 1).Value = ' 1:               search.max.iter                20'
$ws.Cells.Item(196, 1).Value = ' 2:               search.expander                 2'
$ws.Cells.Item(197, 1).Value = ' 3:           search.num.init.exp                 9'
$ws.Cells.Item(198, 1).Value = ' 4:               max.nonconverge              0.01'
$ws.Cells.Item(199, 1).Value = ' 5:                   random.seed             12345'
$ws.Cells.Item(200, 1).Value = ' 6:                output.filestr                NA'
$ws.Cells.Item(201, 1).Value = ' 7:      add.timestamp.to.filestr             FALSE'
$ws.Cells.Item(202, 1).Value = ' 8:           min.obs.date.to.fit                NA'
$ws.Cells.Item(203, 1).Value = ' 9:           max.obs.date.to.fit        2020-04-10'
$ws.Cells.Item(204, 1).Value = '10:               main.iterations             10000'
$ws.Cells.Item(205, 1).Value = '11:         simulation.start.date        2020-01-23'
$ws.Cells.Item(206, 1).Value = '12:        lower.bound.multiplier               0.9'
$ws.Cells.Item(207, 1).Value = '13:        upper.bound.multiplier               1.1'
$ws.Cells.Item(208, 1).Value = '14:            required.in.bounds              0.95'
$ws.Cells.Item(209, 1).Value = '15:                 show.progress              TRUE'
$ws.Cells.Item(210, 1).Value = '16:  plot.observed.data.long.term             FALSE'
$ws.Cells.Item(211, 1).Value = '17: plot.observed.data.short.term              TRUE'
$ws.Cells.Item(212, 1).Value = '18:             lower.bound.label Confirmed COVID19'
$ws.Cells.Item(213, 1).Value = '19:             upper.bound.label  Probable COVID19'
$ws.Cells.Item(214, 1).Value = ''
$ws.Cells.Item(215, 1).Value = '$extras$time.of.run'
$ws.Cells.Item(216, 1).Value = '[1] "2020-05-07 13:56:09"'
$ws.Cells.Item(217, 1).Value = ''
$ws.Cells.Item(218, 1).Value = '$extras$LEMMA.version'
$ws.Cells.Item(219, 1).Value = '     version '
$ws.Cells.Item(220, 1).Value = '"0.3.0.9004" '
$ws.Cells.Item(221, 1).Value = ''
$ws.Cells.Item(222, 1).Value = ''
$ws.Cells.Item(223, 1).Value = '$hosp.bounds'
$ws.Cells.Item(224, 1).Value = '          date    lower     upper'
$ws.Cells.Item(225, 1).Value = ' 1: 2020-04-01 51.88321  59.66569'
$ws.Cells.Item(226, 1).Value = ' 2: 2020-04-02 56.86563  65.39548'
$ws.Cells.Item(227, 1).Value = ' 3: 2020-04-03 61.68362  70.93617'
$ws.Cells.Item(228, 1).Value = ' 4: 2020-04-04 66.26620  76.20613'
$ws.Cells.Item(229, 1).Value = ' 5: 2020-04-05 70.60557  81.19641'
$ws.Cells.Item(230, 1).Value = ' 6: 2020-04-06 74.69862  85.90341'
$ws.Cells.Item(231, 1).Value = ' 7: 2020-04-07 78.82360  90.64714'
$ws.Cells.Item(232, 1).Value = ' 8: 2020-04-08 82.50392  94.87951'
$ws.Cells.Item(233, 1).Value = ' 9: 2020-04-09 85.20577  97.98664'
$ws.Cells.Item(234, 1).Value = '10: 2020-04-10 87.19809 100.27781'
$ws.Cells.Item(235, 1).Value = '11: 2020-04-11 88.66911 101.96948'
$ws.Cells.Item(236, 1).Value = '12: 2020-04-12 89.34864 102.75094'
$ws.Cells.Item(237, 1).Value = '13: 2020-04-13 88.25265 101.49055'
$ws.Cells.Item(238, 1).Value = '14: 2020-04-14 86.63990  99.63588'
$ws.Cells.Item(239, 1).Value = '15: 2020-04-15 85.23594  98.02133'
$ws.Cells.Item(240, 1).Value = '16: 2020-04-17 84.00665  96.60765'
$ws.Cells.Item(241, 1).Value = '17: 2020-04-18 83.83743  96.41305'
$ws.Cells.Item(242, 1).Value = '18: 2020-04-19 83.86506  96.44481'
$ws.Cells.Item(243, 1).Value = '19: 2020-04-20 83.89666  96.48116'
$ws.Cells.Item(244, 1).Value = '20: 2020-04-21 83.85270  96.43061'
$ws.Cells.Item(245, 1).Value = '21: 2020-04-22 83.14715  95.61923'
$ws.Cells.Item(246, 1).Value = '22: 2020-04-23 82.40143  94.76164'
$ws.Cells.Item(247, 1).Value = '23: 2020-04-24 82.17643  94.50289'
$ws.Cells.Item(248, 1).Value = '24: 2020-04-25 81.80076  94.07088'
$ws.Cells.Item(249, 1).Value = '25: 2020-04-26 81.04030  93.19635'
$ws.Cells.Item(250, 1).Value = '26: 2020-04-27 80.07035  92.08090'
$ws.Cells.Item(251, 1).Value = '27: 2020-04-28 78.38229  90.13963'
$ws.Cells.Item(252, 1).Value = '28: 2020-04-29 76.64307  88.13953'
$ws.Cells.Item(253, 1).Value = '29: 2020-04-30 75.31391  86.61100'
$ws.Cells.Item(254, 1).Value = '30: 2020-05-01 74.15291  85.27584'
$ws.Cells.Item(255, 1).Value = '31: 2020-05-02 73.32236  84.32072'
$ws.Cells.Item(256, 1).Value = '32: 2020-05-03 72.73400  83.64410'
$ws.Cells.Item(257, 1).Value = '          date    lower     upper'
$ws.Cells.Item(258, 1).Value = ''
$ws.Cells.Item(259, 1).Value = '$internal.args'
$ws.Cells.Item(260, 1).Value = '$internal.args$search.max.iter'
$ws.Cells.Item(261, 1).Value = '[1] 20'
$ws.Cells.Item(262, 1).Value = ''
$ws.Cells.Item(263, 1).Value = '$internal.args$search.expander'
$ws.Cells.Item(264, 1).Value = '[1] 2'
$ws.Cells.Item(265, 1).Value = ''
$ws.Cells.Item(266, 1).Value = '$internal.args$search.num.init.exp'
$ws.Cells.Item(267, 1).Value = '[1] 9'
$ws.Cells.Item(268, 1).Value = ''
$ws.Cells.Item(269, 1).Value = '$internal.args$max.nonconverge'
$ws.Cells.Item(270, 1).Value = '[1] 0.01'
$ws.Cells.Item(271, 1).Value = ''
$ws.Cells.Item(272, 1).Value = '$internal.args$random.seed'
$ws.Cells.Item(273, 1).Value = '[1] 12345'
$ws.Cells.Item(274, 1).Value = ''
$ws.Cells.Item(275, 1).Value = '$internal.args$output.filestr'
$ws.Cells.Item(276, 1).Value = '[1] "Alameda-May4-v8 output"'
$ws.Cells.Item(277, 1).Value = ''
$ws.Cells.Item(278, 1).Value = '$internal.args$add.timestamp.to.filestr'
$ws.Cells.Item(279, 1).Value = '[1] FALSE'
$ws.Cells.Item(280, 1).Value = ''
$ws.Cells.Item(281, 1).Value = '$internal.args$min.obs.date.to.fit'
$ws.Cells.Item(282, 1).Value = '[1] NA'
$ws.Cells.Item(283, 1).Value = ''
$ws.Cells.Item(284, 1).Value = '$internal.args$max.obs.date.to.fit'
$ws.Cells.Item(285, 1).Value = '[1] "2020-04-10"'
$ws.Cells.Item(286, 1).Value = ''
$ws.Cells.Item(287, 1).Value = '$internal.args$main.iterations'
$ws.Cells.Item(288, 1).Value = '[1] 10000'
$ws.Cells.Item(289, 1).Value = ''
$ws.Cells.Item(290, 1).Value = '$internal.args$simulation.start.date'
$ws.Cells.Item(291, 1).Value = '[1] "2020-01-23"'
$ws.Cells.Item(292, 1).Value = ''
$ws.Cells.Item(293, 1).Value = '$internal.args$lower.bound.multiplier'
$ws.Cells.Item(294, 1).Value = '[1] 0.9'
$ws.Cells.Item(295, 1).Value = ''
$ws.Cells.Item(296, 1).Value = '$internal.args$upper.bound.multiplier'
$ws.Cells.Item(297, 1).Value = '[1] 1.1'
$ws.Cells.Item(298, 1).Value = ''
$ws.Cells.Item(299, 1).Value = '$internal.args$required.in.bounds'
$ws.Cells.Item(300, 1).Value = '[1] 0.95'
$ws.Cells.Item(301, 1).Value = ''
$ws.Cells.Item(302, 1).Value = '$internal.args$show.progress'
$ws.Cells.Item(303, 1).Value = '[1] TRUE'
$ws.Cells.Item(304, 1).Value = ''
$ws.Cells.Item(305, 1).Value = '$internal.args$plot.observed.data.long.term'
$ws.Cells.Item(306, 1).Value = '[1] FALSE'
$ws.Cells.Item(307, 1).Value = ''
$ws.Cells.Item(308, 1).Value = '$internal.args$plot.observed.data.short.term'
$ws.Cells.Item(309, 1).Value = '[1] TRUE'
$ws.Cells.Item(310, 1).Value = ''
$ws.Cells.Item(311, 1).Value = '$internal.args$lower.bound.label'
$ws.Cells.Item(312, 1).Value = '[1] "Confirmed COVID19"'
$ws.Cells.Item(313, 1).Value = ''
$ws.Cells.Item(314, 1).Value = '$internal.args$upper.bound.label'
$ws.Cells.Item(315, 1).Value = '[1] "Probable COVID19"'
$ws.Cells.Item(316, 1).Value = ''
$ws.Cells.Item(317, 1).Value = ''
$ws.Cells.Item(318, 1).Value = '$model.inputs'
$ws.Cells.Item(319, 1).Value = '$model.inputs$total.population'
$ws.Cells.Item(320, 1).Value = '[1] 1671000'
$ws.Cells.Item(321, 1).Value = ''
$ws.Cells.Item(322, 1).Value = '$model.inputs$start.display.date'
$ws.Cells.Item(323, 1).Value = '[1] "2020-03-15"'
$ws.Cells.Item(324, 1).Value = ''
$ws.Cells.Item(325, 1).Value = '$model.inputs$end.date'
$ws.Cells.Item(326, 1).Value = '[1] "2020-07-01"'
$ws.Cells.Item(327, 1).Value = ''
$ws.Cells.Item(328, 1).Value = ''
$ws.Cells.Item(329, 1).Value = '$observed.data'
$ws.Cells.Item(330, 1).Value = '          date     hosp'
$ws.Cells.Item(331, 1).Value = ' 1: 2020-04-01 55.77445'
$ws.Cells.Item(332, 1).Value = ' 2: 2020-04-02 61.13056'
$ws.Cells.Item(333, 1).Value = ' 3: 2020-04-03 66.30990'
$ws.Cells.Item(334, 1).Value = ' 4: 2020-04-04 71.23616'
$ws.Cells.Item(335, 1).Value = ' 5: 2020-04-05 75.90099'
$ws.Cells.Item(336, 1).Value = ' 6: 2020-04-06 80.30101'
$ws.Cells.Item(337, 1).Value = ' 7: 2020-04-07 84.73537'
$ws.Cells.Item(338, 1).Value = ' 8: 2020-04-08 88.69171'
$ws.Cells.Item(339, 1).Value = ' 9: 2020-04-09 91.59621'
$ws.Cells.Item(340, 1).Value = '10: 2020-04-10 93.73795'
$ws.Cells.Item(341, 1).Value = ''
$ws.Cells.Item(342, 1).Value = '$upp.params'
$ws.Cells.Item(343, 1).Value = '   r0.initial latent.period illness.length.given.nonhosp prop.hospitalized'
$ws.Cells.Item(344, 1).Value = '1:          4             3                            5              0.04'
$ws.Cells.Item(345, 1).Value = '   hosp.length.of.stay prop.icu prop.vent intervention1.date'
$ws.Cells.Item(346, 1).Value = '1:                  14     0.45      0.85         2020-03-07'
$ws.Cells.Item(347, 1).Value = '   intervention1.multiplier intervention1.smooth.days intervention2.date'
$ws.Cells.Item(348, 1).Value = '1:                     0.45                         7         2020-03-17'
$ws.Cells.Item(349, 1).Value = '   intervention2.multiplier intervention2.smooth.days intervention3.date'
$ws.Cells.Item(350, 1).Value = '1:                     0.45                         7         2020-06-01'
$ws.Cells.Item(351, 1).Value = '   intervention3.multiplier intervention3.smooth.days'
$ws.Cells.Item(352, 1).Value = '1:                      1.5                         7'
$ws.Cells.Item(353, 1).Value = '   patients.in.hosp.are.infectious use.hosp.rate exposed.to.hospital'
$ws.Cells.Item(354, 1).Value = '1:                            TRUE          TRUE                   9'
$ws.Cells.Item(355, 1).Value = ''
